# Scheduled market-data refresh: update currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) across the leve-profit tables on each job sheet.
# Values come from the upstream price-feed snapshot for this run.
$wb = $excel.ActiveWorkbook

$ws = $wb.Sheets.Item("ALC")
$ws.Range("H51").Value = 3999.5
$ws.Range("I51").Value = 3999.5
$ws.Range("K51").Value = 3999.5
$ws.Range("M51").Value = -3515.5
$ws.Range("H58").Value = 1628.5454
$ws.Range("I58").Value = 989.6
$ws.Range("J58").Value = 2161
$ws.Range("K58").Value = 2968.8
$ws.Range("L58").Value = 6483
$ws.Range("M58").Value = -2818.8
$ws.Range("N58").Value = -6783
$ws.Range("H100").Value = 2987.375
$ws.Range("J100").Value = 3033.3333
$ws.Range("L100").Value = 3033.3333
$ws.Range("N100").Value = -4115.3333
$ws.Range("H137").Value = 1965.619
$ws.Range("I137").Value = 1721.8235
$ws.Range("J137").Value = 3001.75
$ws.Range("K137").Value = 5165.470499999999
$ws.Range("L137").Value = 9005.25
$ws.Range("M137").Value = -2615.470499999999
$ws.Range("N137").Value = -14105.25

$ws = $wb.Sheets.Item("ARM")
$ws.Range("H2").Value = 998.7692
$ws.Range("I2").Value = 1098.4445
$ws.Range("J2").Value = 774.5
$ws.Range("K2").Value = 1098.4445
$ws.Range("L2").Value = 774.5
$ws.Range("M2").Value = -985.4445000000001
$ws.Range("N2").Value = -1000.5
$ws.Range("H5").Value = 105.111115
$ws.Range("I5").Value = 105.111115
$ws.Range("K5").Value = 105.111115
$ws.Range("M5").Value = 6.888885000000002
$ws.Range("H110").Value = 697.6667
$ws.Range("I110").Value = 547.75
$ws.Range("K110").Value = 547.75
$ws.Range("M110").Value = 1497.25
$ws.Range("H116").Value = 998.7692
$ws.Range("I116").Value = 1098.4445
$ws.Range("J116").Value = 774.5
$ws.Range("K116").Value = 1098.4445
$ws.Range("L116").Value = 774.5
$ws.Range("M116").Value = 1195.5555
$ws.Range("N116").Value = -5362.5
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("M122").ClearContents()

$ws = $wb.Sheets.Item("BSM")
$ws.Range("H3").Value = 998.7692
$ws.Range("I3").Value = 1098.4445
$ws.Range("J3").Value = 774.5
$ws.Range("K3").Value = 1098.4445
$ws.Range("L3").Value = 774.5
$ws.Range("M3").Value = -984.4445000000001
$ws.Range("N3").Value = -1002.5
$ws.Range("H4").Value = 105.111115
$ws.Range("I4").Value = 105.111115
$ws.Range("K4").Value = 105.111115
$ws.Range("M4").Value = 9.888885000000002
$ws.Range("H105").Value = 3893.5
$ws.Range("I105").Value = 3735.5715
$ws.Range("K105").Value = 3735.5715
$ws.Range("M105").Value = -1988.5715
$ws.Range("H120").Value = 29000
$ws.Range("J120").Value = 29000
$ws.Range("L120").Value = 29000
$ws.Range("N120").Value = -38676

$ws = $wb.Sheets.Item("CRP")
$ws.Range("H7").Value = 53.42857
$ws.Range("I7").Value = 45.666668
$ws.Range("K7").Value = 45.666668
$ws.Range("M7").Value = 67.333332
$ws.Range("H105").Value = 550
$ws.Range("J105").Value = 0
$ws.Range("L105").Value = 0
$ws.Range("N105").ClearContents()

$ws = $wb.Sheets.Item("CUL")
$ws.Range("H141").Value = 2126.6667
$ws.Range("I141").Value = 2126.6667
$ws.Range("K141").Value = 6380.000100000001
$ws.Range("M141").Value = -1200.000100000001

$ws = $wb.Sheets.Item("GSM")
$ws.Range("H70").Value = 4250
$ws.Range("I70").Value = 2750
$ws.Range("J70").Value = 5000
$ws.Range("K70").Value = 2750
$ws.Range("L70").Value = 5000
$ws.Range("M70").Value = -2480
$ws.Range("N70").Value = -5540
$ws.Range("H73").Value = 4250
$ws.Range("I73").Value = 2750
$ws.Range("J73").Value = 5000
$ws.Range("K73").Value = 2750
$ws.Range("L73").Value = 5000
$ws.Range("M73").Value = -1814
$ws.Range("N73").Value = -6872
$ws.Range("H122").Value = 3491.1428
$ws.Range("I122").Value = 3461.375
$ws.Range("J122").Value = 3530.8333
$ws.Range("K122").Value = 10384.125
$ws.Range("L122").Value = 10592.4999
$ws.Range("M122").Value = -7934.125
$ws.Range("N122").Value = -15492.4999

$ws = $wb.Sheets.Item("LTW")
$ws.Range("H22").Value = 1287
$ws.Range("I22").Value = 1374.5
$ws.Range("J22").Value = 1243.25
$ws.Range("K22").Value = 1374.5
$ws.Range("L22").Value = 1243.25
$ws.Range("M22").Value = -1079.5
$ws.Range("N22").Value = -1833.25
$ws.Range("H27").Value = 1287
$ws.Range("I27").Value = 1374.5
$ws.Range("J27").Value = 1243.25
$ws.Range("K27").Value = 1374.5
$ws.Range("L27").Value = 1243.25
$ws.Range("M27").Value = -1267.5
$ws.Range("N27").Value = -1457.25
$ws.Range("H80").Value = 11000
$ws.Range("J80").Value = 11000
$ws.Range("L80").Value = 11000
$ws.Range("N80").Value = -13246
$ws.Range("H82").Value = 862.375
$ws.Range("I82").Value = 799.8
$ws.Range("K82").Value = 799.8
$ws.Range("M82").Value = -438.8
$ws.Range("H83").Value = 11000
$ws.Range("J83").Value = 11000
$ws.Range("L83").Value = 33000
$ws.Range("N83").Value = -44232
$ws.Range("H85").Value = 862.375
$ws.Range("I85").Value = 799.8
$ws.Range("K85").Value = 799.8
$ws.Range("M85").Value = 448.2
$ws.Range("H127").Value = 77498
$ws.Range("J127").Value = 77498
$ws.Range("L127").Value = 77498
$ws.Range("N127").Value = -87418

$ws = $wb.Sheets.Item("WVR")
$ws.Range("H14").Value = 3440.8
$ws.Range("I14").Value = 2002
$ws.Range("J14").Value = 4400
$ws.Range("K14").Value = 2002
$ws.Range("L14").Value = 4400
$ws.Range("M14").Value = -1834
$ws.Range("N14").Value = -4736
$ws.Range("H81").Value = 2637.5833
$ws.Range("I81").Value = 2786.4546
$ws.Range("K81").Value = 5572.9092
$ws.Range("M81").Value = -4511.9092
$ws.Range("H84").Value = 2637.5833
$ws.Range("I84").Value = 2786.4546
$ws.Range("K84").Value = 27864.546
$ws.Range("M84").Value = -22560.546
$ws.Range("H113").Value = 374
$ws.Range("I113").Value = 417.16666
$ws.Range("J113").Value = 244.5
$ws.Range("K113").Value = 1251.49998
$ws.Range("L113").Value = 733.5
$ws.Range("M113").Value = 918.5000199999999
$ws.Range("N113").Value = -5073.5
$ws.Range("H122").Value = 3843.923
$ws.Range("I122").Value = 3171.75
$ws.Range("J122").Value = 4919.4
$ws.Range("K122").Value = 9515.25
$ws.Range("L122").Value = 14758.2
$ws.Range("M122").Value = -7065.25
$ws.Range("N122").Value = -19658.2
